# Auto-generated script implementing the 2022-Q3 sheet insertion
$wb = $excel.ActiveWorkbook

# ---- Step 1: Update the '总计' (Total) summary sheet ----
$total = $wb.Worksheets.Item(1)

# Read existing data rows (rows 2-8) before overwriting them
$existB = @()
$existC = @()
$existD = @()
for ($r = 2; $r -le 8; $r++) {
    $existB += $total.Range("B$r").Value()
    $existC += $total.Range("C$r").Value()
    $existD += $total.Range("D$r").Value()
}

# Shift rows 2-8 down to 3-9, preserving the A-column style by copying the cell
for ($r = 8; $r -ge 2; $r--) {
    $total.Range("A$r").Copy($total.Range("A" + ($r + 1)))
}
for ($i = 0; $i -lt 7; $i++) {
    $r = $i + 3
    $total.Range("B$r").Value = $existB[$i]
    $total.Range("C$r").Value = $existC[$i]
    $total.Range("D$r").Value = $existD[$i]
}

# Write new row 2 for 2022-Q3
$total.Range("A2").Value = 0
$total.Range("B2").Value = '2022-Q3'
$total.Range("C2").Value = 43
$total.Range("D2").Value = 24.53

# ---- Step 2: Insert the new '2022-Q3' worksheet right after '总计' ----
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = '2022-Q3'

# Header row styled like the other detail sheets (bold/border style copied from 总计!B1)
$headerCols = @('B','C','D','E','F','G','H')
$headerText = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
for ($i = 0; $i -lt 7; $i++) {
    $col = $headerCols[$i]
    $total.Range("B1").Copy($newSheet.Range($col + "1"))
    $newSheet.Range($col + "1").Value = $headerText[$i]
}

# Data rows 2-44 (43 funds). Columns B-G are stored as text, column H as a number,
# matching the source formatting, and column A is a 0-based index styled like 总计!A2.
$dataCols = @('B','C','D','E','F','G')
# row 2: 512880
$total.Range("A2").Copy($newSheet.Range("A2"))
$newSheet.Range("A2").Value = 0
$rowVals = @('512880', '国泰中证全指证券公司ETF', '289.27', '99.95', '2.95', '8.5335')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "2")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H2").Value = 9

# row 3: 512000
$total.Range("A2").Copy($newSheet.Range("A3"))
$newSheet.Range("A3").Value = 1
$rowVals = @('512000', '华宝中证全指证券公司ETF', '215.91', '99.89', '2.90', '6.2614')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "3")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H3").Value = 9

# row 4: 512900
$total.Range("A2").Copy($newSheet.Range("A4"))
$newSheet.Range("A4").Value = 2
$rowVals = @('512900', '南方中证全指证券公司ETF', '78.74', '99.98', '2.91', '2.2913')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "4")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H4").Value = 9

# row 5: 159841
$total.Range("A2").Copy($newSheet.Range("A5"))
$newSheet.Range("A5").Value = 3
$rowVals = @('159841', '天弘中证全指证券公司ETF', '45.70', '99.94', '2.90', '1.3253')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "5")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H5").Value = 9

# row 6: 159993
$total.Range("A2").Copy($newSheet.Range("A6"))
$newSheet.Range("A6").Value = 4
$rowVals = @('159993', '鹏华国证证券龙头ETF', '13.39', '98.06', '6.39', '0.8556')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "6")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H6").Value = 7

# row 7: 161720
$total.Range("A2").Copy($newSheet.Range("A7"))
$newSheet.Range("A7").Value = 5
$rowVals = @('161720', '招商中证全指证券公司指数（LOF）A', '22.28', '94.50', '2.74', '0.6105')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "7")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H7").Value = 9

# row 8: 501016
$total.Range("A2").Copy($newSheet.Range("A8"))
$newSheet.Range("A8").Value = 6
$rowVals = @('501016', '国泰中证申万证券行业指数（LOF）A', '18.10', '93.42', '2.78', '0.5032')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "8")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H8").Value = 9

# row 9: 163113
$total.Range("A2").Copy($newSheet.Range("A9"))
$newSheet.Range("A9").Value = 7
$rowVals = @('163113', '申万菱信中证申万证券行业指数（LOF）A', '16.20', '93.19', '2.73', '0.4423')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "9")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H9").Value = 9

# row 10: 004814
$total.Range("A2").Copy($newSheet.Range("A10"))
$newSheet.Range("A10").Value = 8
$rowVals = @('004814', '中欧红利优享灵活配置混合A', '17.00', '89.59', '2.57', '0.4369')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "10")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H10").Value = 10

# row 11: 161027
$total.Range("A2").Copy($newSheet.Range("A11"))
$newSheet.Range("A11").Value = 9
$rowVals = @('161027', '富国中证全指证券公司指数A', '12.47', '94.32', '2.73', '0.3404')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "11")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H11").Value = 9

# row 12: 502010
$total.Range("A2").Copy($newSheet.Range("A12"))
$newSheet.Range("A12").Value = 10
$rowVals = @('502010', '易方达证券公司指数（LOF）A', '12.08', '94.58', '2.74', '0.3310')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "12")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H12").Value = 9

# row 13: 160633
$total.Range("A2").Copy($newSheet.Range("A13"))
$newSheet.Range("A13").Value = 11
$rowVals = @('160633', '鹏华中证全指证券公司指数（LOF）A', '11.96', '94.06', '2.73', '0.3265')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "13")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H13").Value = 9

# row 14: 501048
$total.Range("A2").Copy($newSheet.Range("A14"))
$newSheet.Range("A14").Value = 12
$rowVals = @('501048', '汇添富中证全指证券公司指数（LOF）C', '9.31', '93.53', '2.73', '0.2542')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "14")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H14").Value = 9

# row 15: 515010
$total.Range("A2").Copy($newSheet.Range("A15"))
$newSheet.Range("A15").Value = 13
$rowVals = @('515010', '华夏中证全指证券公司ETF', '8.47', '99.66', '2.88', '0.2439')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "15")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H15").Value = 9

# row 16: 004815
$total.Range("A2").Copy($newSheet.Range("A16"))
$newSheet.Range("A16").Value = 14
$rowVals = @('004815', '中欧红利优享灵活配置混合C', '8.56', '89.59', '2.57', '0.2200')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "16")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H16").Value = 10

# row 17: 160516
$total.Range("A2").Copy($newSheet.Range("A17"))
$newSheet.Range("A17").Value = 15
$rowVals = @('160516', '博时中证全指证券公司指数', '7.00', '93.76', '2.71', '0.1897')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "17")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H17").Value = 9

# row 18: 012044
$total.Range("A2").Copy($newSheet.Range("A18"))
$newSheet.Range("A18").Value = 16
$rowVals = @('012044', '鹏华中证全指证券公司指数（LOF）C', '5.89', '94.06', '2.73', '0.1608')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "18")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H18").Value = 9

# row 19: 501047
$total.Range("A2").Copy($newSheet.Range("A19"))
$newSheet.Range("A19").Value = 17
$rowVals = @('501047', '汇添富中证全指证券公司指数（LOF）A', '5.58', '93.53', '2.73', '0.1523')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "19")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H19").Value = 9

# row 20: 159842
$total.Range("A2").Copy($newSheet.Range("A20"))
$newSheet.Range("A20").Value = 18
$rowVals = @('159842', '银华中证全指证券公司ETF', '5.01', '98.00', '2.84', '0.1423')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "20")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H20").Value = 9

# row 21: 398041
$total.Range("A2").Copy($newSheet.Range("A21"))
$newSheet.Range("A21").Value = 19
$rowVals = @('398041', '中海量化策略混合', '2.44', '88.10', '5.45', '0.1330')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "21")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H21").Value = 8

# row 22: 515560
$total.Range("A2").Copy($newSheet.Range("A22"))
$newSheet.Range("A22").Value = 20
$rowVals = @('515560', '建信中证全指证券公司ETF', '3.98', '98.63', '2.86', '0.1138')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "22")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H22").Value = 9

# row 23: 502053
$total.Range("A2").Copy($newSheet.Range("A23"))
$newSheet.Range("A23").Value = 21
$rowVals = @('502053', '长盛中证全指证券公司指数（LOF）', '3.93', '93.15', '2.72', '0.1069')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "23")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H23").Value = 9

# row 24: 160419
$total.Range("A2").Copy($newSheet.Range("A24"))
$newSheet.Range("A24").Value = 22
$rowVals = @('160419', '华安中证证券公司A', '3.88', '94.43', '2.72', '0.1055')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "24")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H24").Value = 9

# row 25: 013659
$total.Range("A2").Copy($newSheet.Range("A25"))
$newSheet.Range("A25").Value = 23
$rowVals = @('013659', '中融金融鑫选3个月持有混合A', '1.32', '85.15', '5.15', '0.0680')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "25")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H25").Value = 7

# row 26: 012874
$total.Range("A2").Copy($newSheet.Range("A26"))
$newSheet.Range("A26").Value = 24
$rowVals = @('012874', '易方达证券公司指数（LOF）C', '2.22', '94.58', '2.74', '0.0608')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "26")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H26").Value = 9

# row 27: 512570
$total.Range("A2").Copy($newSheet.Range("A27"))
$newSheet.Range("A27").Value = 25
$rowVals = @('512570', '易方达中证全指证券公司ETF', '1.92', '98.99', '2.87', '0.0551')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "27")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H27").Value = 9

# row 28: 515850
$total.Range("A2").Copy($newSheet.Range("A28"))
$newSheet.Range("A28").Value = 26
$rowVals = @('515850', '富国中证全指证券公司ETF', '1.67', '99.74', '2.77', '0.0463')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "28")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H28").Value = 9

# row 29: 013660
$total.Range("A2").Copy($newSheet.Range("A29"))
$newSheet.Range("A29").Value = 27
$rowVals = @('013660', '中融金融鑫选3个月持有混合C', '0.81', '85.15', '5.15', '0.0417')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "29")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H29").Value = 7

# row 30: 510200
$total.Range("A2").Copy($newSheet.Range("A30"))
$newSheet.Range("A30").Value = 28
$rowVals = @('510200', '汇安上证证券ETF', '0.67', '95.06', '4.47', '0.0299')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "30")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H30").Value = 7

# row 31: 159848
$total.Range("A2").Copy($newSheet.Range("A31"))
$newSheet.Range("A31").Value = 29
$rowVals = @('159848', '国联安中证全指证券公司ETF', '0.96', '96.87', '2.82', '0.0271')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "31")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H31").Value = 9

# row 32: 008116
$total.Range("A2").Copy($newSheet.Range("A32"))
$newSheet.Range("A32").Value = 30
$rowVals = @('008116', '银华沪深股通精选混合', '0.50', '88.46', '4.75', '0.0238')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "32")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H32").Value = 6

# row 33: 516730
$total.Range("A2").Copy($newSheet.Range("A33"))
$newSheet.Range("A33").Value = 31
$rowVals = @('516730', '浦银安盛中证证券公司30ETF', '0.60', '97.43', '3.62', '0.0217')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "33")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H33").Value = 9

# row 34: 013276
$total.Range("A2").Copy($newSheet.Range("A34"))
$newSheet.Range("A34").Value = 32
$rowVals = @('013276', '富国中证全指证券公司指数C', '0.57', '94.32', '2.73', '0.0156')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "34")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H34").Value = 9

# row 35: 090011
$total.Range("A2").Copy($newSheet.Range("A35"))
$newSheet.Range("A35").Value = 33
$rowVals = @('090011', '大成核心双动力混合', '0.24', '92.56', '5.41', '0.0130')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "35")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H35").Value = 3

# row 36: 013597
$total.Range("A2").Copy($newSheet.Range("A36"))
$newSheet.Range("A36").Value = 34
$rowVals = @('013597', '招商中证全指证券公司指数（LOF）C', '0.39', '94.50', '2.74', '0.0107')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "36")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H36").Value = 9

# row 37: 516200
$total.Range("A2").Copy($newSheet.Range("A37"))
$newSheet.Range("A37").Value = 35
$rowVals = @('516200', '华安中证全指证券公司ETF', '0.32', '97.22', '2.90', '0.0093')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "37")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H37").Value = 9

# row 38: 000417
$total.Range("A2").Copy($newSheet.Range("A38"))
$newSheet.Range("A38").Value = 36
$rowVals = @('000417', '国联安新精选灵活配置混合', '0.48', '37.31', '1.58', '0.0076')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "38")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H38").Value = 10

# row 39: 015859
$total.Range("A2").Copy($newSheet.Range("A39"))
$newSheet.Range("A39").Value = 37
$rowVals = @('015859', '宝盈国证证券龙头指数A', '0.12', '94.13', '6.12', '0.0073')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "39")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H39").Value = 7

# row 40: 014984
$total.Range("A2").Copy($newSheet.Range("A40"))
$newSheet.Range("A40").Value = 38
$rowVals = @('014984', '华安中证证券公司C', '0.18', '94.43', '2.72', '0.0049')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "40")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H40").Value = 9

# row 41: 015860
$total.Range("A2").Copy($newSheet.Range("A41"))
$newSheet.Range("A41").Value = 39
$rowVals = @('015860', '宝盈国证证券龙头指数C', '0.05', '94.13', '6.12', '0.0031')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "41")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H41").Value = 7

# row 42: 015178
$total.Range("A2").Copy($newSheet.Range("A42"))
$newSheet.Range("A42").Value = 40
$rowVals = @('015178', '申万菱信中证申万证券行业指数（LOF）C', '0.08', '93.19', '2.73', '0.0022')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "42")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H42").Value = 9

# row 43: 519117
$total.Range("A2").Copy($newSheet.Range("A43"))
$newSheet.Range("A43").Value = 41
$rowVals = @('519117', '浦银安盛基本面400指数', '0.22', '91.78', '0.81', '0.0018')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "43")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H43").Value = 1

# row 44: 015598
$total.Range("A2").Copy($newSheet.Range("A44"))
$newSheet.Range("A44").Value = 42
$rowVals = @('015598', '国泰中证申万证券行业指数（LOF）C', '0.01', '93.42', '2.78', '0.0003')
for ($j = 0; $j -lt 6; $j++) {
    $col = $dataCols[$j]
    $cell = $newSheet.Range($col + "44")
    $cell.NumberFormat = "@"
    $cell.Value = $rowVals[$j]
}
$newSheet.Range("H44").Value = 9

